# Simulated Wild Card round and logged it
# Update the "R" (road/away? - row 3) target-depth totals on both the
# OFF and DEF sheets to reflect the additional simulated playoff game.

$wb = $excel.ActiveWorkbook

$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 207
$wsOff.Range("C3").Value = 148
$wsOff.Range("D3").Value = 46
$wsOff.Range("E3").Value = 21
$wsOff.Range("F3").Value = 5
$wsOff.Range("G3").Value = 3

$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 261
$wsDef.Range("C3").Value = 170
$wsDef.Range("D3").Value = 67
$wsDef.Range("E3").Value = 27
$wsDef.Range("F3").Value = 7
$wsDef.Range("G3").Value = 3
